# Insert a new row at position 91 (pushes existing rows 91-110 down to 92-111)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new data record
$ws.Cells.Item(91, 1).Value = 2
$ws.Cells.Item(91, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(91, 3).Value = "Coquimbo"
$ws.Cells.Item(91, 4).Value = 45204
$ws.Cells.Item(91, 5).Value = 4
$ws.Cells.Item(91, 6).Value = 100112022
$ws.Cells.Item(91, 7).Value = "Arveja Verde"
$ws.Cells.Item(91, 8).Value = "Sin especificar"
$ws.Cells.Item(91, 9).Value = "Primera"
$ws.Cells.Item(91, 10).Value = 400
$ws.Cells.Item(91, 11).Value = 20000
$ws.Cells.Item(91, 12).Value = 21000
$ws.Cells.Item(91, 13).Value = 20500
$ws.Cells.Item(91, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(91, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(91, 16).Value = 820
$ws.Cells.Item(91, 17).Value = 25
$ws.Cells.Item(91, 18).Value = "Hortaliza"
